$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that followed the title ---
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# --- Step 2: before the final paragraph, insert a new bold "Play Genghis' Reel..." ---
#     paragraph, and change the final paragraph's italic text to the meta description. ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$fullRange = $lastPara.Range
$contentRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Genghis' Reel for Free – Review of World Match's Slot Game</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience the legend of Genghis Khan in this exciting slot game, with Wild and Scatter symbols, and Free Spin bonuses. Play for free now.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$contentRange.InsertXML($xml)
